# Manual Punch Import.xlsx -- header/columns rework
# (StaffId / ApplicationType / SelectPunch / InPunch / OutPunch / Remarks)
# plus a couple of pre-formatted duration cells further out on row 2,
# new workbook theme accent colors, and tidied-up column widths/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "StaffId"
$ws.Range("B1").Value = "ApplicationType"
$ws.Range("C1").Value = "SelectPunch"
$ws.Range("D1").Value = "InPunch"
$ws.Range("E1").Value = "OutPunch"
$ws.Range("F1").Value = "Remarks"

# --- Drop the old sample data row, keep only the formatted duration cells
$ws.Range("B2:F2").ClearContents()
$ws.Range("D2").NumberFormat = "mm:ss.0"
$ws.Range("E2").NumberFormat = "mm:ss.0"
$ws.Range("K2").NumberFormat = "mm:ss.0"
$ws.Range("M2").NumberFormat = "mm:ss.0"

# --- Column widths (tightened to fit the new headers) -----------------
$ws.Columns.Item(1).ColumnWidth = 5.833333333333333
$ws.Columns.Item(2).ColumnWidth = 14.166666666666666
$ws.Columns.Item(3).ColumnWidth = 11.166666666666666
$ws.Columns.Item(4).ColumnWidth = 7.333333333333333
$ws.Columns.Item(5).ColumnWidth = 8.833333333333334
$ws.Columns.Item(6).ColumnWidth = 7.833333333333333

# --- Selection cursor ---------------------------------------------------
$ws.Range("F3").Select() | Out-Null

# --- Workbook theme accent colors (new default Office theme palette) ---
# NOTE: Theme.ThemeColorScheme.Colors(i).RGB takes a BGR-packed long
# (VBA/COM RGB() convention), so bytes are swapped from the target RRGGBB.
$wb.Theme.ThemeColorScheme.Colors(3).RGB  = 0x41280E   # dk2      -> 0E2841
$wb.Theme.ThemeColorScheme.Colors(4).RGB  = 0xE8E8E8   # lt2      -> E8E8E8
$wb.Theme.ThemeColorScheme.Colors(5).RGB  = 0x825F14   # accent1  -> 145F82
$wb.Theme.ThemeColorScheme.Colors(6).RGB  = 0x3173E8   # accent2  -> E87331
$wb.Theme.ThemeColorScheme.Colors(7).RGB  = 0x246C18   # accent3  -> 186C24
$wb.Theme.ThemeColorScheme.Colors(8).RGB  = 0xD59E0F   # accent4  -> 0F9ED5
$wb.Theme.ThemeColorScheme.Colors(9).RGB  = 0x932BA0   # accent5  -> A02B93
$wb.Theme.ThemeColorScheme.Colors(10).RGB = 0x2EA74E   # accent6  -> 4EA72E
$wb.Theme.ThemeColorScheme.Colors(11).RGB = 0x867846   # hlink    -> 467886
$wb.Theme.ThemeColorScheme.Colors(12).RGB = 0x7D6096   # folHlink -> 96607D
